$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B: replace numeric judge numbers with "jN" text labels
for ($r = 2; $r -le 21; $r++) {
    $n = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 2).Value = "j$n"
}

# Update selection to C7
$ws.Range("C7").Select()
